# CRUD Done for Ordering microservice..
#
# Slide 5 ("Open-closed Principle (OCP)") gets its three supporting
# pictures (the diagram images) wrapped into a single group shape, and
# the slide's click-animation that used to target the first picture is
# retargeted to the new group.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# --- touch the title run text (no content change, keeps parity with the
#     author's resave of the title placeholder) -----------------------
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Text = $titleRange.Text

# --- group the three picture shapes (170, 171, 172) -------------------
# They are the 3rd, 4th and 5th shapes in the slide's shape tree.
$pics = $s.Shapes.Range(@(3, 4, 5))
$group = $pics.Group()

# --- retarget the "appear" animation from the first picture to the
#     newly created group -----------------------------------------------
$mainSeq = $s.TimeLine.MainSequence
for ($i = $mainSeq.Count; $i -ge 1; $i--) {
    $effect = $mainSeq.Item($i)
    if ($effect.Shape.Id -eq $group.Id) {
        # skip the effect(s) already bound to the group itself
        continue
    }
    if ($effect.Shape.Name -eq "Google Shape;170;g2315e7a7f5d_0_7") {
        $effect.Delete()
        $null = $mainSeq.AddEffect($group, 10, 0, 1)
    }
}
